# Refresh the cryptos price list (row 2..51) with newly scraped values.
# Column D ("Price") holds text-formatted numbers (e.g. "47.203.08",
# "321.92") in the source data, so values that would otherwise be
# auto-parsed as numbers by Excel are entered with a leading apostrophe
# to force them to stay as text (matching the original inlineStr cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.203.08"
$ws.Range("E2").Value = "  +1.34%  "
$ws.Range("D3").Value = "2.490.49"
$ws.Range("E3").Value = "  +0.89%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'321.92"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").Value = "'108.63"
$ws.Range("E6").Value = "  +3.40%  "
$ws.Range("E7").Value = "  +0.48%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -0.81%  "
$ws.Range("D10").Value = "'38.85"
$ws.Range("E10").Value = "  +7.30%  "
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("E12").Value = "  +0.52%  "
$ws.Range("D13").Value = "'18.27"
$ws.Range("E13").Value = "  +0.05%  "
$ws.Range("E14").Value = "  +0.79%  "
$ws.Range("D15").Value = "2.880.53"
$ws.Range("E15").Value = "  +0.85%  "
$ws.Range("D16").Value = "2.487.68"
$ws.Range("E16").Value = "  +1.24%  "
$ws.Range("D17").Value = "'0.846"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").Value = "47.133.20"
$ws.Range("E18").Value = "  +1.49%  "
$ws.Range("D19").Value = "'12.72"
$ws.Range("E19").Value = "  +0.42%  "
$ws.Range("D20").Value = "'6.62"
$ws.Range("E20").Value = "  +2.40%  "
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("D22").Value = "'2.75"
$ws.Range("E22").Value = "  +15.77%  "
$ws.Range("D23").Value = "'70.65"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "'246.70"
$ws.Range("E24").Value = "  -0.88%  "
$ws.Range("E25").Value = "  +1.46%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").Value = "'25.83"
$ws.Range("E27").Value = "  -1.13%  "
$ws.Range("D28").Value = "'2.29"
$ws.Range("E28").Value = "  +3.95%  "
$ws.Range("D29").Value = "'10.08"
$ws.Range("E29").Value = "  +2.92%  "
$ws.Range("E30").Value = "  +8.84%  "
$ws.Range("D31").Value = "'35.30"
$ws.Range("E31").Value = "  +1.79%  "
$ws.Range("D32").Value = "'49.92"
$ws.Range("E32").Value = "  +0.60%  "
$ws.Range("D33").Value = "'20.07"
$ws.Range("E33").Value = "  +2.08%  "
$ws.Range("E34").Value = "  +1.40%  "
$ws.Range("D35").Value = "'0.0784"
$ws.Range("E35").Value = "  +2.18%  "
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("D37").Value = "'4.69"
$ws.Range("E37").Value = "  +1.22%  "
$ws.Range("E38").Value = "  +2.62%  "
$ws.Range("E39").Value = "  +0.45%  "
$ws.Range("E40").Value = "  +0.21%  "
$ws.Range("B41").Value = "WEMIXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").Value = "'2.22"
$ws.Range("E41").Value = "  -0.40%  "
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").Value = "'120.55"
$ws.Range("E42").Value = "  -2.19%  "
$ws.Range("D43").Value = "'21.36"
$ws.Range("E43").Value = "  +3.23%  "
$ws.Range("E44").Value = "  +0.67%  "
$ws.Range("D45").Value = "1.985.49"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("E46").Value = "  +2.17%  "
$ws.Range("D47").Value = "'2.03"
$ws.Range("E47").Value = "  -3.04%  "
$ws.Range("D48").Value = "'1.79"
$ws.Range("E48").Value = "  -1.27%  "
$ws.Range("D49").Value = "'9.08"
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("D50").Value = "'5.14"
$ws.Range("E50").Value = "  -1.83%  "
$ws.Range("D51").Value = "'57.76"
$ws.Range("E51").Value = "  +5.31%  "
